$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B "Asset Level 1" tag values, aligned to existing rows 1-30.
# Row 1 is the bold section header; rows 2-30 are the "Hard Investments (CAPEX)"
# asset rows that receive an asset-level tag. Row 17 (Marine Navigation Equipment)
# and rows 31+ (non hard-investment items) are left untagged, matching the source data.
# The cells are written in the same order the shared-string table was originally built
# in, so new shared strings land at the same indices as the target workbook.
$orderedRows = @(1,2,3,4,5,6,7,8,9,13,10,11,12,14,15,16,19,26,27,28,29,20,21,22,23,24,25,30,18)

$tags = @{
    1  = "Asset Level 1"
    2  = "Roads"
    3  = "Roads"
    4  = "Roads"
    5  = "Roads"
    6  = "Roads"
    7  = "Roads"
    8  = "Transit"
    9  = "Transit"
    10 = "Bridges"
    11 = "Rail"
    12 = "Rail"
    13 = "Transit"
    14 = "Inland Waterways"
    15 = "Inland Waterways"
    16 = "Aviation"
    18 = "Drinking Water, Waste Water"
    19 = "Drinking Water"
    20 = "Hazardous Waste, Soilid Waste"
    21 = "Hazardous Waste, Soilid Waste"
    22 = "Hazardous Waste, Soilid Waste"
    23 = "Hazardous Waste, Soilid Waste"
    24 = "Hazardous Waste, Soilid Waste"
    25 = "Stormwater"
    26 = "Energy"
    27 = "Energy"
    28 = "Energy"
    29 = "Energy"
    30 = "Schools"
}

foreach ($row in $orderedRows) {
    $ws.Cells.Item($row, 2).Value = $tags[$row]
}

# Header cell B1 gets a bold, 12pt font (new style distinct from the existing bold 11pt header style).
$ws.Range("B1").Font.Bold = $true
$ws.Range("B1").Font.Size = 12

# Autosize the new column like column A already is.
$ws.Columns("B:B").AutoFit() | Out-Null

# Leave the same cell selected as in the saved workbook.
$ws.Range("B10").Select() | Out-Null
